$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 10
$ws.Range("B3").Value = 5.5

# Update the active cell selection on the sheet
$ws.Range("C4").Select()
